# Applies the "first great refactoring finished" edit described by the diff.
#
# Two textual changes inside word/document.xml:
#
# 1) In the "По Симпсону" paragraph, a missing space is inserted between
#    "КДО ЛЖ" and the opening "{{" of the simpson_end_diastolic_volume
#    placeholder.
#
# 2) In the "Правый желудочек" paragraph, everything from the closing
#    ")" of "(N<4,2 см)" through the end of the paragraph (the old
#    "средний", "Толщина передней стенки", and "TAPSE" placeholders/text)
#    is collapsed into three new "_full" placeholders immediately
#    followed by a period.

$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# Only touch the "{{ simpson_end_diastolic_volume" span (it starts exactly
# at a run boundary and shares identical run formatting with the runs that
# follow it), so the preceding "КДО" / " ЛЖ" runs are left untouched and
# the merged run keeps the correct (lang="en-US") character formatting.
$r1 = $d.Content.Find.Execute(
    "{{ simpson_end_diastolic_volume",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " {{ simpson_end_diastolic_volume",
    2)
Write-Output "Change1 (КДО ЛЖ space): $r1"

# --- Change 2 -----------------------------------------------------------
$oldTail = "N<4,2 см), средний  {{ right_ventricle_medium }} см (N< 3,5 см). Толщина передней стенки ПЖ:   {{ right_ventricle_wall_thickness }}  см (N<0,5 см). TAPSE:   {{ tapse }}  см (N>=1,7 см)."
$newTail = "N<4,2 см){{ right_ventricle_medium_full }}{{ right_ventricle_wall_thickness_full }}{{ tapse_full }}."

$r2 = $d.Content.Find.Execute(
    $oldTail,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $newTail,
    2)
Write-Output "Change2 (right ventricle tail collapse): $r2"
